# Regen save_data: column G ("K" = strikeouts, formerly "Strike#") is
# recalculated for each game row (rows 2-12). Write the new K values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 5
    9  = 2
    10 = 2
    11 = 1
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
